# This script reworks the "Docentes responsáveis" / syllabus / evaluation
# block of the LOQ4239 sheet: labels shift up while a couple of rows end up
# re-using stray leftover values, the long syllabus/bibliography texts are
# dropped, and the sheet shrinks from 23 to 21 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: objectives value is replaced by the first teacher's line ------
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# --- Row 13: becomes "Programa resumido:" label; B/C reuse "01/01/2021" ---
# (copy values only from B8/C8 so the text stays a literal string and isn't
#  re-interpreted as a date, while keeping B13/C13's existing style)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: becomes "Short syllabus:" label only, old B/C values removed -
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: becomes "Programa:" label, values reuse the teacher's line ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C15").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: becomes "Syllabus:" label only -------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: becomes "Avaliação:" label only, reverts to default height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).AutoFit()

# --- Row 18: becomes "Método:" label with the second teacher's line -------
# B18 is a brand new cell in column B, whose column has overlapping style
# rules, so fix its format explicitly before writing the text.
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: becomes "Critério:" label with the "Aulas expositivas..." text
$ws.Range("A19").Value = "Critério:"
$ws.Range("B10").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: becomes "Norma de recuperação:" label with the "Média..." text
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: becomes "Bibliografia:" label with the "NF = (MF + PR)/2" text
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Rows.Item(21).RowHeight = 120

# --- Rows 22-23 (old "Norma de recuperação:" and "Bibliografia:" rows,
#     including the long bibliography text) are no longer needed ----------
$ws.Range("A22:A23").EntireRow.Delete()
